# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (H) and "Correspond Handback DateTime" (K)
# columns for the first data row (15fd3d9f-...md) on the zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-26 00:47:35"
$wsZhCn.Range("K2").Value = "2016-08-26 00:47:51"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-26 00:47:39"
$wsDeDe.Range("K2").Value = "2016-08-26 00:47:57"
